# "Generate Report for Archive"
#
# The localization status report is regenerated: every row still showing
# the old "Ready for handoff" status is now "In Translation" (Overview!E:F
# and the per-locale zh-cn / de-de sheets' Status column, rows 2-4). Once
# the text is updated, re-autofit those status columns so their width
# reflects the new (shorter) text, matching what the reporting tool does
# whenever it regenerates the workbook.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# Overview sheet: zh-cn / de-de status columns are E and F, data rows 2-4.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F4").Value = $newStatus
$wsOverview.Columns.Item(5).AutoFit()
$wsOverview.Columns.Item(6).AutoFit()

# Per-locale detail sheets: Status is column C, data rows 2-4.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2:C4").Value = $newStatus
$wsZhCn.Columns.Item(3).AutoFit()

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2:C4").Value = $newStatus
$wsDeDe.Columns.Item(3).AutoFit()
